# Applies the commit "Update on 07 Feb 2025 at 17:42":
#   1. The (single, empty) paragraph in the document body gets the text
#      "Test 10." typed into it.
#   2. The section's page size gains an explicit portrait orientation
#      (w:pgSz/@w:orient="portrait").

$d = $word.ActiveDocument

# --- 1. Type "Test 10." into the document's only (empty) paragraph ---------
$para = $d.Paragraphs(1)
$para.Range.InsertAfter("Test 10.")

# --- 2. Force the page orientation to (explicit) portrait ------------------
$d.PageSetup.Orientation = 0   # wdOrientPortrait
